$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (ALC)
$ws.Range("H4").Value = 4352.5
$ws.Range("I4").Value = 2990.4285
$ws.Range("J4").Value = 7530.6665
$ws.Range("K4").Value = 2990.4285
$ws.Range("L4").Value = 7530.6665
$ws.Range("M4").Value = -2876.4285
$ws.Range("N4").Value = -7758.6665

# Row 12 (ALC)
$ws.Range("H12").Value = 116.333336
$ws.Range("I12").Value = 99.5
$ws.Range("K12").Value = 99.5
$ws.Range("M12").Value = 70.5

# Row 55 (ALC)
$ws.Range("H55").Value = 111.44444
$ws.Range("I55").Value = 114
$ws.Range("J55").Value = 106.333336
$ws.Range("K55").Value = 114
$ws.Range("L55").Value = 106.333336
$ws.Range("M55").Value = 100
$ws.Range("N55").Value = -534.333336

# Row 70 (ALC)
$ws.Range("H70").Value = 17816.143
$ws.Range("I70").Value = 3833.3333
$ws.Range("J70").Value = 28303.25
$ws.Range("K70").Value = 11499.9999
$ws.Range("L70").Value = 84909.75
$ws.Range("M70").Value = -11229.9999
$ws.Range("N70").Value = -85449.75

# Row 73 (ALC)
$ws.Range("H73").Value = 17816.143
$ws.Range("I73").Value = 3833.3333
$ws.Range("J73").Value = 28303.25
$ws.Range("K73").Value = 11499.9999
$ws.Range("L73").Value = 84909.75
$ws.Range("M73").Value = -10563.9999
$ws.Range("N73").Value = -86781.75

# Row 121 (ALC)
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 10 (ARM)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 32 (ARM)
$ws.Range("H32").Value = 657.5625
$ws.Range("I32").Value = 535.4643
$ws.Range("K32").Value = 535.4643
$ws.Range("M32").Value = -248.4643

# Row 61 (ARM)
$ws.Range("H61").Value = 1638.6
$ws.Range("I61").Value = 1638.6
$ws.Range("K61").Value = 1638.6
$ws.Range("M61").Value = -1426.6

# Row 74 (ARM)
$ws.Range("H74").Value = 920.2
$ws.Range("I74").Value = 920.2
$ws.Range("K74").Value = 920.2
$ws.Range("M74").Value = -46.20000000000005

# Row 77 (ARM)
$ws.Range("H77").Value = 920.2
$ws.Range("I77").Value = 920.2
$ws.Range("K77").Value = 4601
$ws.Range("M77").Value = -233

# Row 88 (ARM)
$ws.Range("H88").Value = 4668.3335
$ws.Range("I88").Value = 4002.5
$ws.Range("J88").Value = 6000
$ws.Range("K88").Value = 4002.5
$ws.Range("L88").Value = 6000
$ws.Range("M88").Value = -3596.5
$ws.Range("N88").Value = -6812

# Row 91 (ARM)
$ws.Range("H91").Value = 4668.3335
$ws.Range("I91").Value = 4002.5
$ws.Range("J91").Value = 6000
$ws.Range("K91").Value = 4002.5
$ws.Range("L91").Value = 6000
$ws.Range("M91").Value = -2598.5
$ws.Range("N91").Value = -8808

# Row 132 (ARM)
$ws.Range("H132").Value = 1368
$ws.Range("I132").Value = 1368
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4104
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1574
$ws.Range("N132").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 1638.6
$ws.Range("I136").Value = 1638.6
$ws.Range("K136").Value = 4915.799999999999
$ws.Range("M136").Value = -2365.799999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 11 (CRP)
$ws.Range("H11").Value = 150
$ws.Range("J11").Value = 150
$ws.Range("L11").Value = 150
$ws.Range("N11").Value = -430

# Row 31 (CRP)
$ws.Range("H31").Value = 2841.85
$ws.Range("I31").Value = 1142.8334
$ws.Range("J31").Value = 3570
$ws.Range("K31").Value = 1142.8334
$ws.Range("L31").Value = 3570
$ws.Range("M31").Value = -847.8334
$ws.Range("N31").Value = -4160

# Row 34 (CRP)
$ws.Range("H34").Value = 2841.85
$ws.Range("I34").Value = 1142.8334
$ws.Range("J34").Value = 3570
$ws.Range("K34").Value = 1142.8334
$ws.Range("L34").Value = 3570
$ws.Range("M34").Value = -940.8334
$ws.Range("N34").Value = -3974

# Row 111 (CRP)
$ws.Range("H111").Value = 74751
$ws.Range("J111").Value = 74751
$ws.Range("L111").Value = 74751
$ws.Range("N111").Value = -82931

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (CUL)
$ws.Range("H4").Value = 263403.25
$ws.Range("I4").Value = 263403.25
$ws.Range("K4").Value = 790209.75
$ws.Range("M4").Value = -790097.75

# Row 56 (CUL)
$ws.Range("H56").Value = 20000
$ws.Range("I56").Value = 20000
$ws.Range("K56").Value = 20000
$ws.Range("M56").Value = -19470

# Row 57 (CUL)
$ws.Range("H57").Value = 1006
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 1006
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 3018
$ws.Range("N57").Value = -4136
$ws.Range("M57").ClearContents()

# Row 104 (CUL)
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

# Row 113 (CUL)
$ws.Range("H113").Value = 567.8570999999999
$ws.Range("I113").Value = 307.33334
$ws.Range("J113").Value = 763.25
$ws.Range("K113").Value = 922.0000200000001
$ws.Range("L113").Value = 2289.75
$ws.Range("M113").Value = 1247.99998
$ws.Range("N113").Value = -6629.75

$ws = $wb.Worksheets.Item("GSM")
# Row 92 (GSM)
$ws.Range("H92").Value = 12242
$ws.Range("J92").Value = 13890.4
$ws.Range("L92").Value = 13890.4
$ws.Range("N92").Value = -17634.4

# Row 96 (GSM)
$ws.Range("H96").Value = 10261
$ws.Range("J96").Value = 10261
$ws.Range("L96").Value = 10261
$ws.Range("N96").Value = -15753

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("I16").Value = 125002100
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 125002100
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -125001930
$ws.Range("N16").ClearContents()

# Row 40 (LTW)
$ws.Range("H40").Value = 9882.412
$ws.Range("I40").Value = 10031.3125
$ws.Range("K40").Value = 10031.3125
$ws.Range("M40").Value = -9895.3125

# Row 46 (LTW)
$ws.Range("H46").Value = 2471.2856
$ws.Range("I46").Value = 1579.8
$ws.Range("J46").Value = 4700
$ws.Range("K46").Value = 1579.8
$ws.Range("L46").Value = 4700
$ws.Range("M46").Value = -1391.8
$ws.Range("N46").Value = -5076

# Row 55 (LTW)
$ws.Range("H55").Value = 2741.923
$ws.Range("I55").Value = 1162.4286
$ws.Range("K55").Value = 1162.4286
$ws.Range("M55").Value = -989.4286

$ws = $wb.Worksheets.Item("WVR")
# Row 21 (WVR)
$ws.Range("H21").Value = 7500000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 35 (WVR)
$ws.Range("H35").Value = 7500000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 96 (WVR)
$ws.Range("H96").Value = 4900
$ws.Range("I96").Value = 3557.1428
$ws.Range("J96").Value = 7250
$ws.Range("K96").Value = 3557.1428
$ws.Range("L96").Value = 7250
$ws.Range("M96").Value = -2184.1428
$ws.Range("N96").Value = -9996

# Row 136 (WVR)
$ws.Range("H136").Value = 4873.25
$ws.Range("I136").Value = 3812.75
$ws.Range("J136").Value = 6994.25
$ws.Range("K136").Value = 11438.25
$ws.Range("L136").Value = 20982.75
$ws.Range("M136").Value = -8888.25
$ws.Range("N136").Value = -26082.75
